$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert two new rows right after the existing data row (row 16),
#    shifting the blank rows and the signature block down by two rows.
# ------------------------------------------------------------------
$ws.Range("17:18").Insert(-4121) | Out-Null   # xlShiftDown = -4121

# Copy the formatting of the original data row (row 16) onto the two
# freshly inserted rows so they look like the rest of the table.
$ws.Range("B16:J16").Copy() | Out-Null
$ws.Range("B17:J17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B18:J18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Update the header figures.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 48214      # VALOR MORA
$ws.Range("C13").Value = 2          # Cant. Trabajadores
$ws.Range("F13").Value = 2          # Cant. Periodos

# ------------------------------------------------------------------
# 3) Write the worker rows.
#    Row 16 / Row 17 -> new worker (DANIEL BRIEVA MEZA), two periods.
#    Row 18           -> previous worker (ANDRES FELIPE PONCE MORALES).
# ------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1050950344"
$ws.Range("D16").Value = "DANIEL BRIEVA MEZA"
$ws.Range("E16").Value = "2207"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 1000000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1050950344"
$ws.Range("D17").Value = "DANIEL BRIEVA MEZA"
$ws.Range("E17").Value = "2207"
$ws.Range("F17").Value = 6667
$ws.Range("G17").Value = 1000000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1050970745"
$ws.Range("D18").Value = "ANDRES FELIPE PONCE MORALES"
$ws.Range("E18").Value = "2309"
$ws.Range("F18").Value = 1547
$ws.Range("G18").Value = 1300000

Write-Host "Edit applied"
